# Logged Week 17 data: add RB G.Dortch row to the "RB" sheet and make that
# sheet the active tab (it was previously "TE").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RB")

# Switch focus to the RB sheet (moves tabSelected off of TE and onto RB,
# and updates the workbook's activeTab).
$ws.Activate()

# New player row: name in column A, stat columns B:J all zero.
$ws.Range("A6").Value = "G.Dortch"
$ws.Range("B6:J6").Value = 0

# Leave the selection on the newly added name cell.
$ws.Range("A6").Select()
